$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 12.402677
$ws.Range("H2").Value = 37.20803100000001
$ws.Range("I2").Value = 0.4952943482020729
$ws.Range("J2").Value = 0.495294348202073
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.720721333333333
$ws.Range("N2").Value = 23.162164
$ws.Range("O2").Value = 0.1914196631940246
$ws.Range("P2").Value = 0.1914196631940246
$ws.Range("Q2").Value = 95.75761290434268
$ws.Range("R2").Value = 861.8185161390841
$ws.Range("S2").Value = 0.09480907731474476
$ws.Range("T2").Value = 0.09480907731474476

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 12.402677
$ws.Range("H3").Value = 37.20803100000001
$ws.Range("I3").Value = 0.4952943482020729
$ws.Range("J3").Value = 0.495294348202073
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 18.940215
$ws.Range("N3").Value = 56.820645
$ws.Range("O3").Value = 0.4695843069053151
$ws.Range("P3").Value = 0.4695843069053151
$ws.Range("Q3").Value = 234.909368955555
$ws.Range("R3").Value = 2114.184320599995
$ws.Range("S3").Value = 0.2325824532145902
$ws.Range("T3").Value = 0.2325824532145902

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 12.402677
$ws.Range("H4").Value = 37.20803100000001
$ws.Range("I4").Value = 0.4952943482020729
$ws.Range("J4").Value = 0.495294348202073
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.673067
$ws.Range("N4").Value = 41.019201
$ws.Range("O4").Value = 0.3389960299006603
$ws.Range("P4").Value = 0.3389960299006603
$ws.Range("Q4").Value = 169.582633600359
$ws.Range("R4").Value = 1526.243702403231
$ws.Range("S4").Value = 0.167902817672738
$ws.Range("T4").Value = 0.167902817672738

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.161818666666666
$ws.Range("H5").Value = 15.485456
$ws.Range("I5").Value = 0.2061344991927113
$ws.Range("J5").Value = 0.2061344991927113
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.720721333333333
$ws.Range("N5").Value = 23.162164
$ws.Range("O5").Value = 0.1914196631940246
$ws.Range("P5").Value = 0.1914196631940246
$ws.Range("Q5").Value = 39.85296349853155
$ws.Range("R5").Value = 358.676671486784
$ws.Range("S5").Value = 0.03945819640813774
$ws.Range("T5").Value = 0.03945819640813775

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.161818666666666
$ws.Range("H6").Value = 15.485456
$ws.Range("I6").Value = 0.2061344991927113
$ws.Range("J6").Value = 0.2061344991927113
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 18.940215
$ws.Range("N6").Value = 56.820645
$ws.Range("O6").Value = 0.4695843069053151
$ws.Range("P6").Value = 0.4695843069053151
$ws.Range("Q6").Value = 97.76595533767998
$ws.Range("R6").Value = 879.89359803912
$ws.Range("S6").Value = 0.09679752593268355
$ws.Range("T6").Value = 0.09679752593268358

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.161818666666666
$ws.Range("H7").Value = 15.485456
$ws.Range("I7").Value = 0.2061344991927113
$ws.Range("J7").Value = 0.2061344991927113
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 13.673067
$ws.Range("N7").Value = 41.019201
$ws.Range("O7").Value = 0.3389960299006603
$ws.Range("P7").Value = 0.3389960299006603
$ws.Range("Q7").Value = 70.57789247118399
$ws.Range("R7").Value = 635.2010322406559
$ws.Range("S7").Value = 0.06987877685188999
$ws.Range("T7").Value = 0.06987877685189001

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.073119
$ws.Range("H8").Value = 9.219357
$ws.Range("I8").Value = 0.1227233823836907
$ws.Range("J8").Value = 0.1227233823836907
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.720721333333333
$ws.Range("N8").Value = 23.162164
$ws.Range("O8").Value = 0.1914196631940246
$ws.Range("P8").Value = 0.1914196631940246
$ws.Range("Q8").Value = 23.726695423172
$ws.Range("R8").Value = 213.540258808548
$ws.Range("S8").Value = 0.02349166852191757
$ws.Range("T8").Value = 0.02349166852191757

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.073119
$ws.Range("H9").Value = 9.219357
$ws.Range("I9").Value = 0.1227233823836907
$ws.Range("J9").Value = 0.1227233823836907
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 18.940215
$ws.Range("N9").Value = 56.820645
$ws.Range("O9").Value = 0.4695843069053151
$ws.Range("P9").Value = 0.4695843069053151
$ws.Range("Q9").Value = 58.205534580585
$ws.Range("R9").Value = 523.8498112252651
$ws.Range("S9").Value = 0.05762897445772135
$ws.Range("T9").Value = 0.05762897445772136

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.073119
$ws.Range("H10").Value = 9.219357
$ws.Range("I10").Value = 0.1227233823836907
$ws.Range("J10").Value = 0.1227233823836907
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.673067
$ws.Range("N10").Value = 41.019201
$ws.Range("O10").Value = 0.3389960299006603
$ws.Range("P10").Value = 0.3389960299006603
$ws.Range("Q10").Value = 42.01896198597299
$ws.Range("R10").Value = 378.170657873757
$ws.Range("S10").Value = 0.04160273940405178
$ws.Range("T10").Value = 0.04160273940405179

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.403408
$ws.Range("H11").Value = 13.210224
$ws.Range("I11").Value = 0.175847770221525
$ws.Range("J11").Value = 0.175847770221525
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 7.720721333333333
$ws.Range("N11").Value = 23.162164
$ws.Range("O11").Value = 0.1914196631940246
$ws.Range("P11").Value = 0.1914196631940246
$ws.Range("Q11").Value = 33.99748608497066
$ws.Range("R11").Value = 305.9773747647359
$ws.Range("S11").Value = 0.03366072094922455
$ws.Range("T11").Value = 0.03366072094922455

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.403408
$ws.Range("H12").Value = 13.210224
$ws.Range("I12").Value = 0.175847770221525
$ws.Range("J12").Value = 0.175847770221525
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 18.940215
$ws.Range("N12").Value = 56.820645
$ws.Range("O12").Value = 0.4695843069053151
$ws.Range("P12").Value = 0.4695843069053151
$ws.Range("Q12").Value = 83.40149425271998
$ws.Range("R12").Value = 750.61344827448
$ws.Range("S12").Value = 0.08257535330031993
$ws.Range("T12").Value = 0.08257535330031994

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.403408
$ws.Range("H13").Value = 13.210224
$ws.Range("I13").Value = 0.175847770221525
$ws.Range("J13").Value = 0.175847770221525
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 13.673067
$ws.Range("N13").Value = 41.019201
$ws.Range("O13").Value = 0.3389960299006603
$ws.Range("P13").Value = 0.3389960299006603
$ws.Range("Q13").Value = 60.20809261233599
$ws.Range("R13").Value = 541.8728335110239
$ws.Range("S13").Value = 0.05961169597198053
$ws.Range("T13").Value = 0.05961169597198054
